# Generate Report for Handoff
# The file "8cd6338a-92dc-4ecf-a8b3-c0329448181e.md" has been handed off for
# localization, so its status moves from "In Translation" to
# "Ready for handoff", its priority flips from "ht" to "mt", and its
# "Latest Handoff Datetime" / "Latest HO Xliff Generate Date" timestamps are
# refreshed to reflect the new handoff.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for 8cd6338a-....md (row 3) ---
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E3").Value = "Ready for handoff"   # zh-cn status
$ovw.Range("F3").Value = "Ready for handoff"   # de-de status
$ovw.Range("G3").Value = "2016-08-22 18:14:16" # Latest HO Xliff Generate Date

# --- zh-cn sheet: row for 8cd6338a-....md (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"          # Status
$zhcn.Range("E3").Value = "mt"                         # Priority
$zhcn.Range("H3").Value = "2016-08-22 18:14:11"        # Latest Handoff Datetime

# --- de-de sheet: row for 8cd6338a-....md (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"          # Status
$dede.Range("E3").Value = "mt"                         # Priority
$dede.Range("H3").Value = "2016-08-22 18:14:16"        # Latest Handoff Datetime

# --- Column width tweaks (status columns widen to fit "Ready for handoff") ---
$ovw.Columns.Item(5).ColumnWidth = 16.25
$ovw.Columns.Item(6).ColumnWidth = 16.25
$zhcn.Columns.Item(3).ColumnWidth = 16.25
$dede.Columns.Item(3).ColumnWidth = 16.25
